$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 33; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "First$i"
    $ws.Cells.Item($row, 2).Value = "Last$i"
}
